# CHE_sto_pumpedhydro.xlsx - add an independent "capacity_to_activity" entry
# (Preparation for transport: CAP2ACT is now entity dependent)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the existing "co2_factor" row (row 10), which
# shifts every following row down by one and carries their formatting along.
$ws.Rows.Item(10).Insert()

# The freshly inserted row inherits the style of the row above on one cell;
# strip it so the new row matches a plain, unstyled data row.
$ws.Range("C10").ClearFormats()

# Populate the new "capacity_to_activity" parameter row.
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "sto_elec_pumpedhydro"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# Grow the autofilter / filter-database range by one row to keep covering
# the (still mostly blank) tail of the sheet, as before the edit.
$ws.AutoFilterMode = $false
$ws.Range("A5:L574").AutoFilter()
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$5:`$L`$574"

# Restore the selection to where the editor left off.
$ws.Activate()
$ws.Range("B10").Select()
